$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Dlk1-Notch3 LR-pair results (NATMI re-run per Dr Hou advice):
# full cross product of ECs/FAPs/sCs sending clusters x ECs/FAPs/M2/sCs target clusters.
$data = @(
    @("ECs", "Dlk1", "Notch3", "ECs", 1.0, 0.3333333333333333, 0.028147, 0.084441, 0.0007347845853173872, 0.000734784585317387, 3.0, 1.0, 17.32492333333333, 51.97477, 0.1811369435820938, 0.1811369435820938, 0.4876446170633334, 4.38880155357, 0.0001330966339756278, 0.0001330966339756277),
    @("ECs", "Dlk1", "Notch3", "FAPs", 1.0, 0.3333333333333333, 0.028147, 0.084441, 0.0007347845853173872, 0.000734784585317387, 3.0, 1.0, 6.590479666666667, 19.771439, 0.06890531753540824, 0.06890531753540824, 0.1855022311776667, 1.669520080599, 0.00005063056517141783, 0.00005063056517141782),
    @("ECs", "Dlk1", "Notch3", "M2", 1.0, 0.3333333333333333, 0.028147, 0.084441, 0.0007347845853173872, 0.000734784585317387, 2.0, 0.6666666666666666, 0.111195, 0.333585, 0.001162574982531578, 0.001162574982531578, 0.003129805665, 0.028168250985, 0.000000854242176439834, 0.0000008542421764398338),
    @("ECs", "Dlk1", "Notch3", "sCs", 1.0, 0.3333333333333333, 0.028147, 0.084441, 0.0007347845853173872, 0.000734784585317387, 3.0, 1.0, 71.61884566666667, 214.856537, 0.7487951638999665, 0.7487951638999664, 2.015855648979667, 18.142700840817, 0.0005502031439939018, 0.0005502031439939016),
    @("FAPs", "Dlk1", "Notch3", "ECs", 3.0, 1.0, 36.97491766666666, 110.924753, 0.9652396186039789, 0.9652396186039789, 3.0, 1.0, 17.32492333333333, 51.97477, 0.1811369435820938, 0.1811369435820938, 640.5876138313122, 5765.28852448181, 0.1748405543382707, 0.1748405543382707),
    @("FAPs", "Dlk1", "Notch3", "FAPs", 3.0, 1.0, 36.97491766666666, 110.924753, 0.9652396186039789, 0.9652396186039789, 3.0, 1.0, 6.590479666666667, 19.771439, 0.06890531753540824, 0.06890531753540824, 243.6824430588408, 2193.141987529567, 0.06651014241766351, 0.06651014241766351),
    @("FAPs", "Dlk1", "Notch3", "M2", 3.0, 1.0, 36.97491766666666, 110.924753, 0.9652396186039789, 0.9652396186039789, 2.0, 0.6666666666666666, 0.111195, 0.333585, 0.001162574982531578, 0.001162574982531578, 4.111425969945, 37.002833729505, 0.001122163432737308, 0.001122163432737308),
    @("FAPs", "Dlk1", "Notch3", "sCs", 3.0, 1.0, 36.97491766666666, 110.924753, 0.9652396186039789, 0.9652396186039789, 3.0, 1.0, 71.61884566666667, 214.856537, 0.7487951638999665, 0.7487951638999664, 2648.100921906707, 23832.90829716036, 0.7227667584153075, 0.7227667584153074),
    @("sCs", "Dlk1", "Notch3", "ECs", 3.0, 1.0, 1.303400333333333, 3.910201, 0.03402559681070371, 0.03402559681070371, 3.0, 1.0, 17.32492333333333, 51.97477, 0.1811369435820938, 0.1811369435820938, 22.58131084764111, 203.23179762877, 0.00616329260984751, 0.00616329260984751),
    @("sCs", "Dlk1", "Notch3", "FAPs", 3.0, 1.0, 1.303400333333333, 3.910201, 0.03402559681070371, 0.03402559681070371, 3.0, 1.0, 6.590479666666667, 19.771439, 0.06890531753540824, 0.06890531753540824, 8.590033394359889, 77.310300549239, 0.002344544552573313, 0.002344544552573313),
    @("sCs", "Dlk1", "Notch3", "M2", 3.0, 1.0, 1.303400333333333, 3.910201, 0.03402559681070371, 0.03402559681070371, 2.0, 0.6666666666666666, 0.111195, 0.333585, 0.001162574982531578, 0.001162574982531578, 0.144931600065, 1.304384400585, 0.00003955730761783037, 0.00003955730761783037),
    @("sCs", "Dlk1", "Notch3", "sCs", 3.0, 1.0, 1.303400333333333, 3.910201, 0.03402559681070371, 0.03402559681070371, 3.0, 1.0, 71.61884566666667, 214.856537, 0.7487951638999665, 0.7487951638999664, 93.34802731488189, 840.132245833937, 0.02547820234066506, 0.02547820234066506)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
